$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add_AWB")

# Update the AWB serial numbers in column B (rows 2 and 3),
# consistent with the next values in the sequence. Use a leading
# apostrophe so the numeric-looking text stays stored as text
# (matching the existing quotePrefix text formatting of the column).
$ws.Range("B2").Value = "'9702488"
$ws.Range("B3").Value = "'9702489"

# Move the active selection down one row, as it ends up after the edit.
$ws.Range("C11").Select()
